# Applies the "New crime data collected" update to the CompStat weekly
# report worksheet: refreshes the report header (volume number + the
# covered week's date range) and the weekly/28-day/YTD/2-year crime
# count table for rows 15-30 (Rape .. Hate Crimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Plain numeric write - leaves the cell's existing style/number-format
# untouched (used when a cell stays numeric before and after the edit).
function Set-Num($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Converts a (numeric) cell into a text cell holding $val, and pulls the
# cell style from $styleSrc (a same-row "already text" cell) so the
# resulting style id matches what a genuinely-text cell would carry
# (numFmtId 0 / "right,center" text style) instead of a synthesized
# text-number-format style.
function Set-TextWithStyle($addr, $val, $styleSrc) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $ws.Range($styleSrc).Copy()
    $c.PasteSpecial(-4122)   # xlPasteFormats
}

# Converts a (text) cell into a numeric cell holding $val, and pulls the
# cell style from $styleSrc (a same-row "already numeric" cell) so the
# resulting style id matches a genuinely-numeric cell.
function Set-NumWithStyle($addr, $val, $styleSrc) {
    $c = $ws.Range($addr)
    $c.Value = $val
    $ws.Range($styleSrc).Copy()
    $c.PasteSpecial(-4122)   # xlPasteFormats
}

# ---------------------------------------------------------------------
# Report header: Volume/Number and the covered week's dates.
# ---------------------------------------------------------------------

$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-Num "M15" 25
Set-Num "N15" -44.444444444444

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
Set-Num "C16" 3
Set-Num "D16" 1
Set-Num "E16" 200
Set-Num "F16" 11
Set-Num "G16" 15
Set-Num "H16" -26.666666666666
Set-Num "I16" 37
Set-Num "J16" 44
Set-Num "K16" -15.909090909090
Set-Num "L16" -5.128205128205
Set-Num "M16" 12.121212121212
Set-Num "N16" -83.027522935779

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
Set-Num "C17" 3
Set-Num "D17" 4
Set-Num "E17" -25
Set-Num "F17" 14
Set-Num "G17" 7
Set-Num "H17" 100
Set-Num "I17" 46
Set-Num "J17" 32
Set-Num "K17" 43.75
Set-Num "L17" 119.047619047619
Set-Num "M17" 70.370370370370
Set-Num "N17" 21.052631578947

# ---------------------------------------------------------------------
# Row 18 - Burglary (D18/E18 become "no data" text markers)
# ---------------------------------------------------------------------
Set-Num "C18" 2
Set-TextWithStyle "D18" "0" "D14"
Set-TextWithStyle "E18" "***.*" "E14"
Set-Num "F18" 8
Set-Num "G18" 8
Set-Num "H18" 0
Set-Num "I18" 51
Set-Num "J18" 48
Set-Num "K18" 6.25
Set-Num "L18" 41.666666666666
Set-Num "M18" 18.604651162790
Set-Num "N18" -82.229965156794

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
Set-Num "C19" 23
Set-Num "D19" 14
Set-Num "E19" 64.285714285714
Set-Num "F19" 62
Set-Num "G19" 63
Set-Num "H19" -1.587301587301
Set-Num "I19" 260
Set-Num "J19" 271
Set-Num "K19" -4.059040590405
Set-Num "L19" 23.809523809523
Set-Num "M19" 3.585657370517
Set-Num "N19" -63.838664812239

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-Num "C20" 2
Set-Num "D20" 2
Set-Num "E20" 0
Set-Num "F20" 7
Set-Num "G20" 7
Set-Num "H20" 0
Set-Num "I20" 29
Set-Num "J20" 21
Set-Num "K20" 38.095238095238
Set-Num "L20" 52.631578947368
Set-Num "M20" 190
Set-Num "N20" -93.363844393592

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
Set-Num "C21" 33
Set-Num "D21" 21
Set-Num "E21" 57.142857142857
Set-Num "F21" 102
Set-Num "G21" 100
Set-Num "H21" 2
Set-Num "I21" 429
Set-Num "J21" 421
Set-Num "K21" 1.900237529691
Set-Num "L21" 30.395136778115
Set-Num "M21" 16.260162601626
Set-Num "N21" -74.89760093622

# ---------------------------------------------------------------------
# Row 22 - Transit (D22/E22 flip from "no data" text back to numeric)
# ---------------------------------------------------------------------
Set-Num "C22" 1
Set-NumWithStyle "D22" 1 "D16"
Set-NumWithStyle "E22" 0 "E16"
Set-Num "F22" 2
Set-Num "G22" 3
Set-Num "H22" -33.333333333333
Set-Num "I22" 11
Set-Num "J22" 9
Set-Num "K22" 22.222222222222
Set-Num "L22" 37.5
Set-Num "M22" 10

# ---------------------------------------------------------------------
# Row 23 - Housing (C23 becomes a "no data" text marker)
# ---------------------------------------------------------------------
Set-TextWithStyle "C23" "0" "D26"
Set-Num "D23" 1
Set-Num "E23" -100
Set-Num "F23" 3
Set-Num "G23" 2
Set-Num "H23" 50
Set-Num "I23" 18
Set-Num "J23" 5
Set-Num "K23" 260
Set-Num "L23" 80

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
Set-Num "C24" 9
Set-Num "D24" 25
Set-Num "E24" -64
Set-Num "F24" 85
Set-Num "G24" 104
Set-Num "H24" -18.269230769230
Set-Num "I24" 383
Set-Num "J24" 507
Set-Num "K24" -24.457593688362
Set-Num "L24" -28.411214953271
Set-Num "M24" 7.887323943661

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
Set-Num "C25" 3
Set-Num "D25" 4
Set-Num "E25" -25
Set-Num "F25" 17
Set-Num "G25" 17
Set-Num "H25" 0
Set-Num "I25" 83
Set-Num "J25" 78
Set-Num "K25" 6.410256410256
Set-Num "L25" 38.333333333333
Set-Num "M25" -10.752688172043

# ---------------------------------------------------------------------
# Row 26 - UCR Rape* (C26/F26 flip from "no data" text back to numeric)
# ---------------------------------------------------------------------
Set-NumWithStyle "C26" 1 "I26"
Set-NumWithStyle "F26" 1 "I26"
Set-Num "I26" 8
Set-Num "K26" 60
Set-Num "L26" 33.333333333333

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes (C27 becomes a "no data" text marker)
# ---------------------------------------------------------------------
Set-TextWithStyle "C27" "0" "D26"
Set-Num "D27" 3
Set-Num "E27" -100
Set-Num "G27" 7
Set-Num "H27" -42.857142857142
Set-Num "J27" 18
Set-Num "K27" -27.777777777777
Set-Num "L27" 30

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes (F30/I30 flip from "no data" text back to numeric)
# ---------------------------------------------------------------------
Set-NumWithStyle "F30" 1 "J30"
Set-NumWithStyle "I30" 1 "J30"
Set-Num "K30" -50
Set-Num "L30" -50
